$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.302.26"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.18%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.879.21"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.97%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.20%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.85%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.0000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.19%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4808"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.59%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2882"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.06%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06583"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.48%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.877.73"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.85%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.89"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.24%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07373"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.35%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.200"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.38%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.01"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.73%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6605"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.34%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.262.23"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.23%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.51"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.31%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.000"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.08%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007722"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.71%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.453"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.88%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.143.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.17%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9996"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.32%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "193.35"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.93%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.177"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.08%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.438"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.98%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.74"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.68%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.28"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.26%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.934"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.27%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.444"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.58%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.272"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.38%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09151"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.33%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.048"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.09%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05068"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.10%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7431"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.27%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.139"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.22%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.715"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.39%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01838"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.43%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.636"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.20%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9155"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.21%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.078"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.21%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "106.44"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.05%  "

# Row 42
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4329"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.69%  "

# Row 43
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.884"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.53%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9991"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.34%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.660"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.60%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1350"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.99%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.583"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +10.06%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "65.42"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -10.03%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.921"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.71%  "

# Row 50
$ws.Range("E50").Value = "  -2.79%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05726"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.71%  "
